$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Designator D13 -> D21 (row 12)
# Leading apostrophe keeps the existing quotePrefix-bearing cell style intact
# (matches original author's edit, which only changed the shared-string value).
$ws.Range("A12").Value = "'D21"

# 2. U13 (Option) -> "U13 (Option)  " (trailing spaces) (row 49)
$ws.Range("A49").Value = "'U13 (Option)  "

# 3. Digikey part number 1528-1438-ND -> LCSC Part # C89297 (row 49)
$ws.Range("C49").Value = "'LCSC Part # C89297"

# 4. Y1 -> "Y1 Buy from Aliexpress, search for 28.375MHz" (row 50)
$ws.Range("A50").Value = "'Y1 Buy from Aliexpress, search for 28.375MHz"

# 5. New row 58: a plain note
$ws.Range("A58").Value = "Note: Y1 buy from here"

# 6. New row 59: hyperlink-styled cell with the Aliexpress URL as text.
# Set the display text first, then attach a (targetless) hyperlink marker so
# the cell picks up the built-in "Hyperlink" style - matches the target,
# which has <hyperlink ref="A59"/> with no r:id (no live target).
$url = "https://www.aliexpress.us/item/3256804924643676.html?spm=a2g0o.productlist.main.1.69ecb531MQJD0Y&algo_pvid=f8189ddb-8fab-4334-927c-212d7a5233a3&algo_exp_id=f8189ddb-8fab-4334-927c-212d7a5233a3-0&pdp_npi=4%40dis%21USD%213.02%211.54%21%21%2122.00%21%21%40210321dc16976695767143414ecc84%2112000031697748677%21sea%21US%21703686525%21&curPageLogUid=fvYFQqMhVbyQ"
$ws.Range("A59").Value = $url
$ws.Hyperlinks.Add($ws.Range("A59"), "")

# Update view / selection state to match target
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A60").Select()

# Page setup scale change
$ws.PageSetup.Zoom = 68
